# Implement the "Assassin" card (and its supporting "The Contract" lord),
# matching the target commit "Implemented the Assassin card".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 41: "Assassin" (MINION) ------------------------------------
# Written left-to-right first so any brand-new shared strings land in
# sharedStrings.xml in the same order as the authoritative edit.
$ws.Range("E41").Value2 = "Assassin"
$ws.Range("F41").Value2 = "MINION"
$ws.Range("G41").Value2 = 3
$ws.Range("H41").Value2 = 1
$ws.Range("I41").Value2 = 4
$ws.Range("J41").Value2 = 'When you summon this unit, search "Assassination"'
$ws.Range("K41").Value2 = "Y"

# --- "The Doctor" (row 38) card effect was reworded ----------------------
$ws.Range("J38").Value2 = "At the end of your turn, restore 1 health to all units on your side of the field."

# --- "Assasination" (row 40) is now implemented, mark it done -----------
$ws.Range("K40").Value2 = "Y"

# --- New row 42: "The Contract" (LORD) -----------------------------------
$ws.Range("E42").Value2 = "The Contract"
$ws.Range("F42").Value2 = "LORD"
$ws.Range("G42").Value2 = 0
$ws.Range("H42").Value2 = 0
$ws.Range("I42").Value2 = 0
$ws.Range("J42").Value2 = "At the end of your opponent's turns, search ""Assassin""."
$ws.Range("K42").Value2 = "N"

# --- Reflect where the editor ended up on the sheet ----------------------
$ws.Range("K40").Select()
